$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was incremented by
# one day (45204 -> 45205) for every data row (rows 2 through 472).
$ws.Range("C2:C472").Value = 45205
